$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.82%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.58%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.142"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05787"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.15%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.41%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.06%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.57%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.30%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1366"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.21%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.07062"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.26%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03265"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'12.65%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.09355"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.001529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.89%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'One"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.0005975"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.03%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006018"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.29%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.67%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.161"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.34%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.60%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03303"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.21%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.50%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.182"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-11.92%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04139"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.78%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.86%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'1.13%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004137"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.84%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'3.35%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E40").Value = "'0.60%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.1070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.10%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'KickToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.003531"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-33.70%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002437"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.06%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009169"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.43%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.58%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05795"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-42.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-22.95%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = "Normal"
